$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the player data (rows 2:205) in descending order by Price (column D),
# using the worksheet's Sort object so the persisted sort state is updated too.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add2($ws.Range("D2:D205"), 0, 2, 0, 0) | Out-Null
$sortObj.SetRange($ws.Range("A1:D205"))
$sortObj.Header = 1
$sortObj.MatchCase = $false
$sortObj.Orientation = 1
$sortObj.Apply()

# Restore the selection to match the post-sort view.
$ws.Range("D127:D205").Select()
